$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 539, shifting existing rows 539-647 down to 540-648.
$ws.Rows.Item(539).EntireRow.Insert()

# Populate the newly inserted row 539 with the new data record.
$ws.Cells.Item(539, 1).Value = 3
$ws.Cells.Item(539, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(539, 3).Value = "Coquimbo"
$ws.Cells.Item(539, 4).Value = 45258
$ws.Cells.Item(539, 5).Value = 5
$ws.Cells.Item(539, 6).Value = 100112012
$ws.Cells.Item(539, 7).Value = "Espinaca"
$ws.Cells.Item(539, 8).Value = "Sin especificar"
$ws.Cells.Item(539, 9).Value = "Primera"
$ws.Cells.Item(539, 10).Value = 80
$ws.Cells.Item(539, 11).Value = 5000
$ws.Cells.Item(539, 12).Value = 5000
$ws.Cells.Item(539, 13).Value = 5000
$ws.Cells.Item(539, 14).Value = "$/docena de atados (3 kilos)"
$ws.Cells.Item(539, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(539, 16).Value = 1667
$ws.Cells.Item(539, 17).Value = 3
$ws.Cells.Item(539, 18).Value = "Hortaliza"
